$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 7-11); this also shrinks
# the sheet's used range / dimension down to A1:K6.
$ws.Rows("7:11").Delete()

# Row 2
$ws.Range("A2").Value = "AGENCE KHATABI"
$ws.Range("B2").Value = "FF"
$ws.Range("C2").Value = "54544646446464646464444464"
$ws.Range("D2").Value = "AGENCE 1"
$ws.Range("E2").Value = "BMCI"
$ws.Range("F2").Value = "Point de vente"
$ws.Range("G2").Value = "389/AOURIR/AV1"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 7500
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 7050

# Row 3
$ws.Range("A3").Value = "AGENCE LAHLOU"
$ws.Range("B3").Value = "1098777"
$ws.Range("C3").Value = "121232435465768778798809"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "CIH"
$ws.Range("F3").Value = "Supervision"
$ws.Range("G3").Value = "001/SUP SUD"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 9000
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 8400

# Row 4
$ws.Range("A4").Value = "BAKKALI MOHAMED"
$ws.Range("B4").Value = "B12346"
$ws.Range("C4").Value = "78017053636372722873881919"
$ws.Range("D4").Value = "HASSAN 2"
$ws.Range("E4").Value = "CIH"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 8000
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 7600

# Row 5
$ws.Range("A5").Value = "AGENCE ESSALAM"
$ws.Range("B5").Value = "19087"
$ws.Range("C5").Value = "671721839230232983487766"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "AWB"
$ws.Range("F5").Value = "Point de vente"
$ws.Range("G5").Value = "604/ERRAHMA"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 9999.99
$ws.Range("J5").Value = 666.66
$ws.Range("K5").Value = 9333.33

# Row 6 - totals row (label cells blank/space, only the MT columns carry totals)
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 34499.99
$ws.Range("J6").Value = 2116.66
$ws.Range("K6").Value = 32383.33
